$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark.
#    Today it sits right after "...Refuerza tu aprendizaje: L" (i.e.
#    collapsed, right before the "os caudillos" run). It needs to move
#    down to the "Palabras clave..." paragraph, landing right after the
#    word "caudillo," (which used to read "independencia").
#    We delete it now and re-add it in its new spot once the text below
#    has been edited (so the character offsets are correct).
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ------------------------------------------------------------------
# 2) Fix the keyword list in the "Palabras clave del recurso..."
#    paragraph:
#      " caudillo, independencia,América Latina"
#    becomes
#      " caudillo,independencia,América Latina"
#    i.e. "caudillo, " (the word, comma and trailing space) is removed
#    from where it was, and "caudillo," is inserted immediately before
#    "independencia" instead.
# ------------------------------------------------------------------
$range = $d.Content
$range.Find.Execute(" caudillo, independencia,América", $true, $false, $false, $false, $false, $true, 1, $false, " caudillo,independencia,América", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark right after "caudillo," (i.e.
#    immediately before "independencia,América") in its new location.
# ------------------------------------------------------------------
$full = $d.Content.Text
$pos = $full.IndexOf(" caudillo,independencia") + " caudillo,".Length
$target = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
